# "Generate Report for Archive"
# The localization status changes from "Ready for handoff" to "In Translation"
# everywhere it's referenced (Overview!E2:F2, zh-cn!C2, de-de!C2). Because the
# shared string is used verbatim across the three sheets, updating the cell
# Value in place also updates every other cell bound to the same shared
# string entry.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# The "Status" columns now hold shorter text, so resize them to fit the new
# content (matches the narrower columns produced by the original report
# generator for Overview!E:F, zh-cn!C and de-de!C).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
